# Apply the "2-Datastore" change:
#  - On the Login sheet, fix the capitalization of the "password" label in B1
#    to "Password" (this re-appends a fresh shared string entry, matching the
#    upstream diff's shared-string table reshuffle).
#  - Update the sheet's last active selection to D3:D4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Login")
$ws.Activate()

$ws.Range("B1").Value = "Password"

$ws.Range("D3:D4").Select()
